$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) The three "spacer" paragraphs that sit right after a table (blank line
#    used to create a gap before the next table) had their near-invisible
#    font size (w:sz/w:szCs = 4, i.e. 2pt) bumped very slightly to 6 (3pt).
#    Font.Size drives w:sz, Font.SizeBi drives w:szCs.
# ---------------------------------------------------------------------------
For ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    If ($p.Range.Font.Size -eq 2) {
        $p.Range.Font.Size = 3
        $p.Range.Font.SizeBi = 3
    }
}

# ---------------------------------------------------------------------------
# 2) The two "Contenudecadre" paragraphs inside the last table's cell no
#    longer force full justification.
# ---------------------------------------------------------------------------
For ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    If ($p.Range.ParagraphFormat.Alignment -eq 3) {
        $p.Range.ParagraphFormat.Alignment = 0
    }
}

# ---------------------------------------------------------------------------
# 3) The last table (the "laureats" frame table) gets a touch less left cell
#    margin and a taller first (only) row, shrinking the blank space below
#    it. Do the table-structure edits last since they re-lay-out the table
#    and can invalidate paragraph handles obtained beforehand.
# ---------------------------------------------------------------------------
$lastTable = $d.Tables.Item($d.Tables.Count)
$lastTable.LeftPadding = 2.1                 # 42 dxa (was 44 dxa / 2.2pt)
$lastTable.Rows.Item(1).Height = 101.65      # 2033 dxa (was 1928 dxa / 96.4pt)

Write-Host "done"
